$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Increase row heights for the data rows (rows 3-16 grow from 18.75 to 19.5,
# and the last data row 17 grows from 19.5 to 20.25).
for ($r = 3; $r -le 16; $r++) {
    $ws.Rows.Item($r).RowHeight = 19.5
}
$ws.Rows.Item(17).RowHeight = 20.25

# Reformat the large empty block below the table (B18:I1000): right align the
# (empty) values and switch the font to a plain black Calibri.
$rng = $ws.Range("B18:I1000")
$rng.HorizontalAlignment = -4152
$rng.Font.Name = "Calibri"
$rng.Font.Color = 0
